# Recurso solar y nubosidad actualizado
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Hoja "Solsticio_Invierno_20Jun" (sheet1)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Solsticio_Invierno_20Jun")

# Ancho de columna D: 8 -> 9
$ws1.Columns.Item(4).ColumnWidth = 8.17

$ws1.Range("D10").Value = 45.42
$ws1.Range("E10").Value = 99.91

$ws1.Range("D11").Value = 464.39
$ws1.Range("E11").Value = 62.81

$ws1.Range("D12").Value = 698.86
$ws1.Range("E12").Value = 49.64

$ws1.Range("D13").Value = 871.8
$ws1.Range("E13").Value = 44.11

$ws1.Range("D14").Value = 976.38
$ws1.Range("E14").Value = 42.54

$ws1.Range("D15").Value = 1009.65
$ws1.Range("E15").Value = 44.11

$ws1.Range("D16").Value = 973.19
$ws1.Range("E16").Value = 49.64

$ws1.Range("D17").Value = 874.43
$ws1.Range("E17").Value = 62.81

$ws1.Range("D18").Value = 744.75
$ws1.Range("E18").Value = 99.91

$ws1.Range("D19").Value = 523.26
$ws1.Range("E19").Value = 200

$ws1.Range("D20").Value = 0
$ws1.Range("E20").Value = -100

# ---------------------------------------------------------------------
# Hoja "Solsticio_Verano_21Dic" (sheet2)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Solsticio_Verano_21Dic")

# Ancho de columna D: 9 -> 8
$ws2.Columns.Item(4).ColumnWidth = 7.17

$ws2.Range("D8").Value = 0
$ws2.Range("E8").Value = -100

$ws2.Range("D9").Value = 85.61
$ws2.Range("E9").Value = -59.81

$ws2.Range("D10").Value = 292.6
$ws2.Range("E10").Value = -36.5

$ws2.Range("D11").Value = 374.14
$ws2.Range("E11").Value = -26.54

$ws2.Range("D12").Value = 594.79
$ws2.Range("E12").Value = -21.56

$ws2.Range("D13").Value = 785.83
$ws2.Range("E13").Value = -19.13

$ws2.Range("D14").Value = 918.45
$ws2.Range("E14").Value = -18.39

$ws2.Range("D15").Value = 919.27
$ws2.Range("E15").Value = -19.13

$ws2.Range("D16").Value = 844.7
$ws2.Range("E16").Value = -21.56

$ws2.Range("D17").Value = 698.24
$ws2.Range("E17").Value = -26.54

$ws2.Range("D18").Value = 487.26
$ws2.Range("E18").Value = -36.5

$ws2.Range("D19").Value = 217.96
$ws2.Range("E19").Value = -59.81

$ws2.Range("D20").Value = 0
$ws2.Range("E20").Value = -100

$ws2.Range("D21").Value = 0
$ws2.Range("E21").Value = -100

# ---------------------------------------------------------------------
# Hoja "Resumen_Estadisticas" (sheet3)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Resumen_Estadisticas")

$ws3.Range("F2").Value = 1009.65
$ws3.Range("G2").Value = 299.26
$ws3.Range("H2").Value = 7.182
$ws3.Range("K2").Value = 59.6
$ws3.Range("L2").Value = 2.682

$ws3.Range("F3").Value = 919.27
$ws3.Range("G3").Value = 259.12
$ws3.Range("H3").Value = 6.219
$ws3.Range("K3").Value = -29.92
$ws3.Range("L3").Value = -2.655

# ---------------------------------------------------------------------
# Hoja "Informacion_General" (sheet4)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Informacion_General")

$ws4.Range("B4").Value = "35°"
$ws4.Range("B8").Value = "2025-07-10 12:37:53"

Write-Host "Recurso solar y nubosidad actualizado"
